# Map107 / Map032 scene update: add English (column C) translations next
# to five existing Japanese (column A) lines of dialogue on rows 26, 27,
# 30, 31 and 32. The Japanese text in column A is left untouched.
#
# Single-quoted here-strings (@' ... '@) are used so the text is taken
# completely literally (no PowerShell interpolation, and the apostrophes
# that occur in the English text don't need escaping). The closing '@
# must start at column 0.
#
# Each write is followed by AutoFit() on that row: entering multi-line
# text through .Value otherwise leaves a stray explicit custom row
# height behind, which the real edit does not have.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$s26 = @'

<Lily>Hey, it feels good, right?
You've got such a little knife♥
'@
$ws.Range("C26").Value = $s26
$ws.Rows.Item(26).AutoFit()

$s27 = @'

<Lily>The only weapon you can hold is the one attached to your crotch♥
There's no way to defeat a succubus, right?
Heehee♥
'@
$ws.Range("C27").Value = $s27
$ws.Rows.Item(27).AutoFit()

$s30 = @'
\ n<
[3]> Heehee ♥ 
Just a little teasing made you react like that?
Can you even call yourself a hunter anymore?
'@
$ws.Range("C30").Value = $s30
$ws.Rows.Item(30).AutoFit()

$s31 = @'

<
[3]>Let me tell you something. 
You can't defeat a succubus with something like that, right? 
Weapons don't work. But you know that already.
'@
$ws.Range("C31").Value = $s31
$ws.Rows.Item(31).AutoFit()

$s32 = @'

<
[3]>I hope you like it. 
Show me more♥
 Little succubus hunter♥
'@
$ws.Range("C32").Value = $s32
$ws.Rows.Item(32).AutoFit()
